$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("January-2021")

# --- Update header date row (I1:AM1): shift from Jan-2020 serials to Jan-2021 serials ---
# I1 previously 43831 (2020-01-01) -> should become 44197 (2021-01-01), and so on through AM1 (31 days)
$startSerial = 44197
for ($i = 0; $i -lt 31; $i++) {
    $col = 9 + $i   # column I = 9 .. AM = 39
    $ws.Cells.Item(1, $col).Value = $startSerial + $i
}

# --- Fill in the newly reported "carry amount" for 31-Jan (column AM) for several retailers ---
$ws.Range("AM5").Value = 1040
$ws.Range("AM14").Value = 5200
$ws.Range("AM18").Value = 5200
$ws.Range("AM25").Value = 3120
$ws.Range("AM56").Value = 5200
$ws.Range("AM61").Value = 2080
$ws.Range("AM66").Value = 5200
$ws.Range("AM71").Value = 3120
$ws.Range("AM81").Value = 1040
$ws.Range("AM96").Value = 2080

# --- Update the frozen-pane view / active selection to reflect where the user left off ---
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Application.ActiveWindow.ScrollColumn = 38
$ws.Range("AL27").Select()
